$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")
$ws.Activate()

# ------------------------------------------------------------------
# Insert 11 new rows (30-40) below the current last row (29) by
# copying the block of rows 19-29 downward. This brings along the
# exact per-cell formatting (fonts/styles) used by those rows so the
# new rows visually match the rest of the table (col A style, col B
# style, etc.) instead of falling back to the worksheet/column
# default style.
# ------------------------------------------------------------------
$ws.Range("A19:G29").Copy()
$ws.Range("A30:G40").Insert(-4121)   # xlShiftDown

# Column C in rows 6-29 mostly uses a different (but visually
# identical) font/style than the one used for the new rows in the
# final workbook. Re-apply the format used by C3 (style actually
# used for column C on the new rows) to column C of the new block.
$ws.Cells.Item(3, 3).Copy()
$ws.Range("C30:C40").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# Fill in the new key values (column C). They are entered in the
# same order the author originally typed them in (so the shared
# string table ends up built in that exact order), and the first two
# rows are swapped afterwards to reach the final on-sheet order.
# ------------------------------------------------------------------
$ws.Cells.Item(30, 3).Value = "IND_08_FL_DATE"
$ws.Cells.Item(31, 3).Value = "IND_08_FL_OVERRIDE"
$ws.Cells.Item(32, 3).Value = "IND_01_FL_OVERRIDE"
$ws.Cells.Item(33, 3).Value = "IND_03_FL_OVERRIDE"
$ws.Cells.Item(34, 3).Value = "IND_04_FL_OVERRIDE"
$ws.Cells.Item(35, 3).Value = "IND_05_FL_OVERRIDE"
$ws.Cells.Item(36, 3).Value = "IND_09_FL_OVERRIDE"
$ws.Cells.Item(37, 3).Value = "IND_12_FL_OVERRIDE"
$ws.Cells.Item(38, 3).Value = "IND_FL_PEGG_BILANCI"
$ws.Cells.Item(39, 3).Value = "IND_FL_DATE6M"
$ws.Cells.Item(40, 3).Value = "INDICATOR_33"

# swap rows 30 & 31 values so the final sheet shows OVERRIDE then DATE
$tmp = $ws.Cells.Item(30, 3).Value2
$ws.Cells.Item(30, 3).Value = $ws.Cells.Item(31, 3).Value2
$ws.Cells.Item(31, 3).Value = $tmp

# Columns A, B, E and F keep the same values on every row of the
# table ("CREATE/MODIFY", "LIB_EWS_IT", "String", "String"); make sure
# that holds for all the newly inserted rows as well.
For ($r = 30; $r -le 40; $r++) {
    $ws.Cells.Item($r, 1).Value = "CREATE/MODIFY"
    $ws.Cells.Item($r, 2).Value = "LIB_EWS_IT"
    $ws.Cells.Item($r, 5).Value = "String"
    $ws.Cells.Item($r, 6).Value = "String"
}

# ------------------------------------------------------------------
# Update the view so the window shows the new rows, matching the
# state the workbook was left in (row 25 scrolled into view, C40
# selected as the active cell).
# ------------------------------------------------------------------
$excel.Goto($ws.Range("A25"), $true)
$ws.Range("C40").Select()
